$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.311.18'
$ws.Range("D3").Value = '1.622.29'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +1.23%  '
$ws.Range("E8").Value = '  +1.64%  '
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.96'
$ws.Range("E10").Value = '  +5.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0815'
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("D12").Value = '1.847.96'
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("D13").Value = '1.626.42'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.520'
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D16").Value = '26.332.25'
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.52'
$ws.Range("E17").Value = '  +4.37%  '
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.10'
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("E21").Value = '  +1.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.37'
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("E24").Value = '  +7.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.39'
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.21'
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("E29").Value = '  +2.15%  '
$ws.Range("E30").Value = '  +11.56%  '
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("E32").Value = '  +2.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.96'
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("E34").Value = '  +2.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("D36").Value = '1.181.53'
$ws.Range("E36").Value = '  +5.31%  '
$ws.Range("E37").Value = '  +1.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.810'
$ws.Range("E38").Value = '  +3.92%  '
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("E43").Value = '  +5.15%  '
$ws.Range("D44").Value = '1.758.89'
$ws.Range("E44").Value = '  +1.76%  '
$ws.Range("E45").Value = '  +1.35%  '
$ws.Range("E46").Value = '  +15.86%  '
$ws.Range("E47").Value = '  +2.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.99'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("E51").Value = '  -0.43%  '
